$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated line counts for several classes/handlers (re-measured loc) ---
$ws.Range("C3").Value  = 375   # engine.entity / BallEntity
$ws.Range("C4").Value  = 149   # BlockEntity
$ws.Range("C5").Value  = 119   # ButtonEntity
$ws.Range("C6").Value  = 493   # StickEntity

$ws.Range("C12").Value = 76    # engine.event / IteratedCollisionEvent

$ws.Range("C14").Value = 113   # factories / BorderFactory
$ws.Range("C15").Value = 193   # ItemFactory

$ws.Range("C17").Value = 137   # handlers / ControllerHandler
$ws.Range("C18").Value = 462   # EntityHandler
$ws.Range("C19").Value = 68    # FileHandler
$ws.Range("C20").Value = 162   # HighscoreHandler
$ws.Range("C21").Value = 52    # ItemHandler
$ws.Range("C22").Value = 122   # LanguageHandler
$ws.Range("C23").Value = 186   # LevelHandler
$ws.Range("C24").Value = 286   # OptionsHandler
$ws.Range("C25").Value = 53    # PlayerHandler
$ws.Range("C26").Value = 246   # SoundHandler

$ws.Range("C29").Value = 171   # parameters / Constatns

$ws.Range("C32").Value = 153   # ui / AboutState
$ws.Range("C33").Value = 128   # Breakout
$ws.Range("C34").Value = 485   # GameplayState
$ws.Range("C35").Value = 270   # HighscoreState
$ws.Range("C37").Value = 280   # OptionsState

# --- New change-log entry: "added target indicator for Stick-BOT" ---
$ws.Range("G11").Value = [DateTime]"2017-09-06"
$ws.Range("I11").Value = 5652

# --- Selection moved from I11 to G11 (also clears the stale topLeftCell) ---
$ws.Range("G11").Select()
